$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 69
$ws.Cells.Item(2, 6).Value = 28
$ws.Cells.Item(2, 7).Value = 0.594
$ws.Cells.Item(2, 9).Value = 36.3
$ws.Cells.Item(2, 10).Value = 78.90000000000001
$ws.Cells.Item(2, 14).Value = 0.363
$ws.Cells.Item(2, 15).Value = 18.5
$ws.Cells.Item(2, 25).Value = 4.3
$ws.Cells.Item(2, 26).Value = 19.4
$ws.Cells.Item(2, 28).Value = 98.40000000000001
$ws.Cells.Item(2, 29).Value = 2.1
$ws.Cells.Item(2, 30).Value = 8
$ws.Cells.Item(2, 32).Value = 11
$ws.Cells.Item(2, 35).Value = 22
$ws.Cells.Item(2, 37).Value = 11
$ws.Cells.Item(2, 40).Value = 18
$ws.Cells.Item(2, 41).Value = 21
$ws.Cells.Item(2, 44).Value = 20
$ws.Cells.Item(2, 50).Value = 16
$ws.Cells.Item(2, 53).Value = 21
$ws.Cells.Item(2, 55).Value = 11
$ws.Cells.Item(2, 58).Value = "'2009-03-21"
# Row 3
$ws.Cells.Item(3, 9).Value = 37.4
$ws.Cells.Item(3, 11).Value = 0.485
$ws.Cells.Item(3, 12).Value = 6.4
$ws.Cells.Item(3, 13).Value = 16.3
$ws.Cells.Item(3, 14).Value = 0.392
$ws.Cells.Item(3, 15).Value = 20
$ws.Cells.Item(3, 16).Value = 26
$ws.Cells.Item(3, 17).Value = 0.769
$ws.Cells.Item(3, 19).Value = 31.8
$ws.Cells.Item(3, 21).Value = 22.6
$ws.Cells.Item(3, 26).Value = 23.5
$ws.Cells.Item(3, 27).Value = 22.6
$ws.Cells.Item(3, 28).Value = 101.1
$ws.Cells.Item(3, 29).Value = 8
$ws.Cells.Item(3, 30).Value = 2
$ws.Cells.Item(3, 42).Value = 7
$ws.Cells.Item(3, 43).Value = 16
$ws.Cells.Item(3, 44).Value = 18
$ws.Cells.Item(3, 45).Value = 6
$ws.Cells.Item(3, 49).Value = 7
$ws.Cells.Item(3, 50).Value = 18
$ws.Cells.Item(3, 52).Value = 29
$ws.Cells.Item(3, 54).Value = 10
$ws.Cells.Item(3, 58).Value = "'2009-03-21"
# Row 4
$ws.Cells.Item(4, 4).Value = 69
$ws.Cells.Item(4, 6).Value = 38
$ws.Cells.Item(4, 7).Value = 0.449
$ws.Cells.Item(4, 11).Value = 0.456
$ws.Cells.Item(4, 12).Value = 6.1
$ws.Cells.Item(4, 14).Value = 0.374
$ws.Cells.Item(4, 15).Value = 17.8
$ws.Cells.Item(4, 17).Value = 0.748
$ws.Cells.Item(4, 18).Value = 10.6
$ws.Cells.Item(4, 20).Value = 39.4
$ws.Cells.Item(4, 21).Value = 21.2
$ws.Cells.Item(4, 25).Value = 5.9
$ws.Cells.Item(4, 26).Value = 21.8
$ws.Cells.Item(4, 28).Value = 93.7
$ws.Cells.Item(4, 29).Value = -0.8
$ws.Cells.Item(4, 30).Value = 8
$ws.Cells.Item(4, 31).Value = 18
$ws.Cells.Item(4, 34).Value = 5
$ws.Cells.Item(4, 40).Value = 10
$ws.Cells.Item(4, 41).Value = 25
$ws.Cells.Item(4, 43).Value = 26
$ws.Cells.Item(4, 44).Value = 19
$ws.Cells.Item(4, 45).Value = 25
$ws.Cells.Item(4, 47).Value = 13
$ws.Cells.Item(4, 49).Value = 20
$ws.Cells.Item(4, 50).Value = 13
$ws.Cells.Item(4, 53).Value = 16
$ws.Cells.Item(4, 55).Value = 17
$ws.Cells.Item(4, 58).Value = "'2009-03-21"
# Row 5
$ws.Cells.Item(5, 5).Value = 32
$ws.Cells.Item(5, 6).Value = 37
$ws.Cells.Item(5, 7).Value = 0.464
$ws.Cells.Item(5, 9).Value = 37.7
$ws.Cells.Item(5, 10).Value = 83.40000000000001
$ws.Cells.Item(5, 15).Value = 19.8
$ws.Cells.Item(5, 16).Value = 25.1
$ws.Cells.Item(5, 17).Value = 0.79
$ws.Cells.Item(5, 18).Value = 12.2
$ws.Cells.Item(5, 19).Value = 30.4
$ws.Cells.Item(5, 21).Value = 20.8
$ws.Cells.Item(5, 23).Value = 7.5
$ws.Cells.Item(5, 26).Value = 21.3
$ws.Cells.Item(5, 27).Value = 20.9
$ws.Cells.Item(5, 29).Value = -1
$ws.Cells.Item(5, 30).Value = 8
$ws.Cells.Item(5, 40).Value = 6
$ws.Cells.Item(5, 41).Value = 10
$ws.Cells.Item(5, 44).Value = 5
$ws.Cells.Item(5, 45).Value = 12
$ws.Cells.Item(5, 47).Value = 15
$ws.Cells.Item(5, 48).Value = 22
$ws.Cells.Item(5, 53).Value = 15
$ws.Cells.Item(5, 54).Value = 9
$ws.Cells.Item(5, 58).Value = "'2009-03-21"
# Row 6
$ws.Cells.Item(6, 5).Value = 55
$ws.Cells.Item(6, 6).Value = 13
$ws.Cells.Item(6, 7).Value = 0.8090000000000001
$ws.Cells.Item(6, 9).Value = 36.7
$ws.Cells.Item(6, 10).Value = 78.40000000000001
$ws.Cells.Item(6, 11).Value = 0.468
$ws.Cells.Item(6, 14).Value = 0.386
$ws.Cells.Item(6, 16).Value = 25
$ws.Cells.Item(6, 17).Value = 0.753
$ws.Cells.Item(6, 20).Value = 41.6
$ws.Cells.Item(6, 21).Value = 19.9
$ws.Cells.Item(6, 22).Value = 13
$ws.Cells.Item(6, 25).Value = 4
$ws.Cells.Item(6, 27).Value = 20.6
$ws.Cells.Item(6, 28).Value = 100.2
$ws.Cells.Item(6, 29).Value = 8.9
$ws.Cells.Item(6, 30).Value = 20
$ws.Cells.Item(6, 35).Value = 14
$ws.Cells.Item(6, 38).Value = 4
$ws.Cells.Item(6, 40).Value = 4
$ws.Cells.Item(6, 41).Value = 17
$ws.Cells.Item(6, 45).Value = 9
$ws.Cells.Item(6, 46).Value = 14
$ws.Cells.Item(6, 47).Value = 26
$ws.Cells.Item(6, 48).Value = 8
$ws.Cells.Item(6, 49).Value = 8
$ws.Cells.Item(6, 53).Value = 19
$ws.Cells.Item(6, 58).Value = "'2009-03-21"
# Row 7
$ws.Cells.Item(7, 30).Value = 2
$ws.Cells.Item(7, 37).Value = 12
$ws.Cells.Item(7, 38).Value = 13
$ws.Cells.Item(7, 44).Value = 14
$ws.Cells.Item(7, 45).Value = 7
$ws.Cells.Item(7, 48).Value = 7
$ws.Cells.Item(7, 54).Value = 10
$ws.Cells.Item(7, 58).Value = "'2009-03-21"
# Row 8
$ws.Cells.Item(8, 30).Value = 2
$ws.Cells.Item(8, 44).Value = 13
$ws.Cells.Item(8, 46).Value = 13
$ws.Cells.Item(8, 49).Value = 2
$ws.Cells.Item(8, 58).Value = "'2009-03-21"
# Row 9
$ws.Cells.Item(9, 30).Value = 20
$ws.Cells.Item(9, 47).Value = 16
$ws.Cells.Item(9, 50).Value = 21
$ws.Cells.Item(9, 53).Value = 28
$ws.Cells.Item(9, 58).Value = "'2009-03-21"
# Row 10
$ws.Cells.Item(10, 30).Value = 8
$ws.Cells.Item(10, 47).Value = 14
$ws.Cells.Item(10, 58).Value = "'2009-03-21"
# Row 11
$ws.Cells.Item(11, 40).Value = 11
$ws.Cells.Item(11, 48).Value = 13
$ws.Cells.Item(11, 50).Value = 26
$ws.Cells.Item(11, 51).Value = 23
$ws.Cells.Item(11, 55).Value = 5
$ws.Cells.Item(11, 58).Value = "'2009-03-21"
# Row 12
$ws.Cells.Item(12, 4).Value = 70
$ws.Cells.Item(12, 5).Value = 28
$ws.Cells.Item(12, 7).Value = 0.4
$ws.Cells.Item(12, 10).Value = 86.09999999999999
$ws.Cells.Item(12, 11).Value = 0.449
$ws.Cells.Item(12, 12).Value = 7.9
$ws.Cells.Item(12, 13).Value = 21.1
$ws.Cells.Item(12, 18).Value = 11.3
$ws.Cells.Item(12, 20).Value = 43.5
$ws.Cells.Item(12, 22).Value = 14.9
$ws.Cells.Item(12, 28).Value = 103.8
$ws.Cells.Item(12, 29).Value = -2.3
$ws.Cells.Item(12, 30).Value = 2
$ws.Cells.Item(12, 33).Value = 22
$ws.Cells.Item(12, 34).Value = 11
$ws.Cells.Item(12, 37).Value = 23
$ws.Cells.Item(12, 41).Value = 20
$ws.Cells.Item(12, 44).Value = 12
$ws.Cells.Item(12, 48).Value = 21
$ws.Cells.Item(12, 49).Value = 22
$ws.Cells.Item(12, 52).Value = 27
$ws.Cells.Item(12, 55).Value = 21
$ws.Cells.Item(12, 58).Value = "'2009-03-21"
# Row 13
$ws.Cells.Item(13, 30).Value = 8
$ws.Cells.Item(13, 32).Value = 28
$ws.Cells.Item(13, 33).Value = 28
$ws.Cells.Item(13, 43).Value = 27
$ws.Cells.Item(13, 44).Value = 15
$ws.Cells.Item(13, 45).Value = 26
$ws.Cells.Item(13, 47).Value = 12
$ws.Cells.Item(13, 49).Value = 21
$ws.Cells.Item(13, 51).Value = 22
$ws.Cells.Item(13, 52).Value = 11
$ws.Cells.Item(13, 58).Value = "'2009-03-21"
# Row 14
$ws.Cells.Item(14, 9).Value = 40.7
$ws.Cells.Item(14, 10).Value = 85.2
$ws.Cells.Item(14, 12).Value = 6.8
$ws.Cells.Item(14, 13).Value = 18.7
$ws.Cells.Item(14, 16).Value = 25.8
$ws.Cells.Item(14, 17).Value = 0.77
$ws.Cells.Item(14, 18).Value = 12.5
$ws.Cells.Item(14, 19).Value = 31.8
$ws.Cells.Item(14, 20).Value = 44.4
$ws.Cells.Item(14, 22).Value = 13.7
$ws.Cells.Item(14, 23).Value = 8.5
$ws.Cells.Item(14, 24).Value = 5.4
$ws.Cells.Item(14, 25).Value = 4.6
$ws.Cells.Item(14, 26).Value = 20.6
$ws.Cells.Item(14, 27).Value = 22.3
$ws.Cells.Item(14, 28).Value = 108.1
$ws.Cells.Item(14, 29).Value = 7.7
$ws.Cells.Item(14, 30).Value = 20
$ws.Cells.Item(14, 38).Value = 14
$ws.Cells.Item(14, 42).Value = 8
$ws.Cells.Item(14, 43).Value = 15
$ws.Cells.Item(14, 45).Value = 5
$ws.Cells.Item(14, 49).Value = 3
$ws.Cells.Item(14, 50).Value = 7
$ws.Cells.Item(14, 51).Value = 13
$ws.Cells.Item(14, 52).Value = 13
$ws.Cells.Item(14, 58).Value = "'2009-03-21"
# Row 15
$ws.Cells.Item(15, 4).Value = 68
$ws.Cells.Item(15, 6).Value = 51
$ws.Cells.Item(15, 7).Value = 0.25
$ws.Cells.Item(15, 9).Value = 34.9
$ws.Cells.Item(15, 10).Value = 77.5
$ws.Cells.Item(15, 11).Value = 0.45
$ws.Cells.Item(15, 12).Value = 4.6
$ws.Cells.Item(15, 13).Value = 13.4
$ws.Cells.Item(15, 14).Value = 0.348
$ws.Cells.Item(15, 16).Value = 25.2
$ws.Cells.Item(15, 18).Value = 10.5
$ws.Cells.Item(15, 19).Value = 28.4
$ws.Cells.Item(15, 27).Value = 21.7
$ws.Cells.Item(15, 28).Value = 93.40000000000001
$ws.Cells.Item(15, 29).Value = -6.3
$ws.Cells.Item(15, 30).Value = 20
$ws.Cells.Item(15, 37).Value = 21
$ws.Cells.Item(15, 44).Value = 21
$ws.Cells.Item(15, 49).Value = 11
$ws.Cells.Item(15, 50).Value = 20
$ws.Cells.Item(15, 58).Value = "'2009-03-21"
# Row 16
$ws.Cells.Item(16, 30).Value = 20
$ws.Cells.Item(16, 51).Value = 5
$ws.Cells.Item(16, 58).Value = "'2009-03-21"
# Row 17
$ws.Cells.Item(17, 9).Value = 36.5
$ws.Cells.Item(17, 10).Value = 82.09999999999999
$ws.Cells.Item(17, 11).Value = 0.445
$ws.Cells.Item(17, 14).Value = 0.36
$ws.Cells.Item(17, 15).Value = 20.2
$ws.Cells.Item(17, 19).Value = 28.9
$ws.Cells.Item(17, 20).Value = 41
$ws.Cells.Item(17, 21).Value = 21.6
$ws.Cells.Item(17, 22).Value = 14.3
$ws.Cells.Item(17, 23).Value = 7.3
$ws.Cells.Item(17, 26).Value = 24.4
$ws.Cells.Item(17, 28).Value = 99.2
$ws.Cells.Item(17, 30).Value = 2
$ws.Cells.Item(17, 31).Value = 18
$ws.Cells.Item(17, 32).Value = 19
$ws.Cells.Item(17, 33).Value = 19
$ws.Cells.Item(17, 35).Value = 17
$ws.Cells.Item(17, 41).Value = 6
$ws.Cells.Item(17, 44).Value = 6
$ws.Cells.Item(17, 48).Value = 14
$ws.Cells.Item(17, 49).Value = 15
$ws.Cells.Item(17, 51).Value = 14
$ws.Cells.Item(17, 55).Value = 18
$ws.Cells.Item(17, 58).Value = "'2009-03-21"
# Row 18
$ws.Cells.Item(18, 30).Value = 8
$ws.Cells.Item(18, 34).Value = 14
$ws.Cells.Item(18, 58).Value = "'2009-03-21"
# Row 19
$ws.Cells.Item(19, 30).Value = 8
$ws.Cells.Item(19, 32).Value = 19
$ws.Cells.Item(19, 38).Value = 3
$ws.Cells.Item(19, 40).Value = 7
$ws.Cells.Item(19, 46).Value = 24
$ws.Cells.Item(19, 47).Value = 25
$ws.Cells.Item(19, 54).Value = 15
$ws.Cells.Item(19, 55).Value = 20
$ws.Cells.Item(19, 58).Value = "'2009-03-21"
# Row 20
$ws.Cells.Item(20, 30).Value = 20
$ws.Cells.Item(20, 31).Value = 8
$ws.Cells.Item(20, 33).Value = 8
$ws.Cells.Item(20, 41).Value = 22
$ws.Cells.Item(20, 46).Value = 23
$ws.Cells.Item(20, 49).Value = 13
$ws.Cells.Item(20, 58).Value = "'2009-03-21"
# Row 21
$ws.Cells.Item(21, 4).Value = 68
$ws.Cells.Item(21, 6).Value = 40
$ws.Cells.Item(21, 7).Value = 0.412
$ws.Cells.Item(21, 12).Value = 10.3
$ws.Cells.Item(21, 13).Value = 28.5
$ws.Cells.Item(21, 14).Value = 0.36
$ws.Cells.Item(21, 16).Value = 23.3
$ws.Cells.Item(21, 17).Value = 0.789
$ws.Cells.Item(21, 18).Value = 11.1
$ws.Cells.Item(21, 19).Value = 31.2
$ws.Cells.Item(21, 21).Value = 21.3
$ws.Cells.Item(21, 22).Value = 14.3
$ws.Cells.Item(21, 24).Value = 2.4
$ws.Cells.Item(21, 25).Value = 5.3
$ws.Cells.Item(21, 26).Value = 20.5
$ws.Cells.Item(21, 27).Value = 19.5
$ws.Cells.Item(21, 29).Value = -2.5
$ws.Cells.Item(21, 30).Value = 20
$ws.Cells.Item(21, 31).Value = 21
$ws.Cells.Item(21, 33).Value = 21
$ws.Cells.Item(21, 34).Value = 19
$ws.Cells.Item(21, 44).Value = 16
$ws.Cells.Item(21, 47).Value = 11
$ws.Cells.Item(21, 49).Value = 13
$ws.Cells.Item(21, 52).Value = 10
$ws.Cells.Item(21, 58).Value = "'2009-03-21"
# Row 22
$ws.Cells.Item(22, 30).Value = 8
$ws.Cells.Item(22, 34).Value = 14
$ws.Cells.Item(22, 37).Value = 25
$ws.Cells.Item(22, 45).Value = 11
$ws.Cells.Item(22, 52).Value = 15
$ws.Cells.Item(22, 58).Value = "'2009-03-21"
# Row 23
$ws.Cells.Item(23, 4).Value = 68
$ws.Cells.Item(23, 5).Value = 50
$ws.Cells.Item(23, 7).Value = 0.735
$ws.Cells.Item(23, 9).Value = 36
$ws.Cells.Item(23, 13).Value = 26.5
$ws.Cells.Item(23, 17).Value = 0.723
$ws.Cells.Item(23, 19).Value = 33.3
$ws.Cells.Item(23, 20).Value = 43.3
$ws.Cells.Item(23, 21).Value = 19.5
$ws.Cells.Item(23, 27).Value = 22.4
$ws.Cells.Item(23, 28).Value = 102.1
$ws.Cells.Item(23, 30).Value = 20
$ws.Cells.Item(23, 40).Value = 3
$ws.Cells.Item(23, 41).Value = 11
$ws.Cells.Item(23, 49).Value = 23
$ws.Cells.Item(23, 50).Value = 8
$ws.Cells.Item(23, 53).Value = 7
$ws.Cells.Item(23, 58).Value = "'2009-03-21"
# Row 24
$ws.Cells.Item(24, 35).Value = 13
$ws.Cells.Item(24, 58).Value = "'2009-03-21"
# Row 25
$ws.Cells.Item(25, 4).Value = 68
$ws.Cells.Item(25, 5).Value = 37
$ws.Cells.Item(25, 7).Value = 0.544
$ws.Cells.Item(25, 9).Value = 40.8
$ws.Cells.Item(25, 10).Value = 80.7
$ws.Cells.Item(25, 11).Value = 0.505
$ws.Cells.Item(25, 12).Value = 6.6
$ws.Cells.Item(25, 14).Value = 0.382
$ws.Cells.Item(25, 15).Value = 20.3
$ws.Cells.Item(25, 16).Value = 26.9
$ws.Cells.Item(25, 19).Value = 30.9
$ws.Cells.Item(25, 20).Value = 41.2
$ws.Cells.Item(25, 21).Value = 22.9
$ws.Cells.Item(25, 22).Value = 15.8
$ws.Cells.Item(25, 23).Value = 7
$ws.Cells.Item(25, 24).Value = 4.9
$ws.Cells.Item(25, 26).Value = 20.8
$ws.Cells.Item(25, 27).Value = 22.6
$ws.Cells.Item(25, 28).Value = 108.5
$ws.Cells.Item(25, 29).Value = 1.9
$ws.Cells.Item(25, 30).Value = 20
$ws.Cells.Item(25, 41).Value = 5
$ws.Cells.Item(25, 45).Value = 10
$ws.Cells.Item(25, 48).Value = 29
$ws.Cells.Item(25, 53).Value = 6
$ws.Cells.Item(25, 55).Value = 12
$ws.Cells.Item(25, 58).Value = "'2009-03-21"
# Row 26
$ws.Cells.Item(26, 5).Value = 43
$ws.Cells.Item(26, 6).Value = 26
$ws.Cells.Item(26, 7).Value = 0.623
$ws.Cells.Item(26, 9).Value = 36.4
$ws.Cells.Item(26, 11).Value = 0.461
$ws.Cells.Item(26, 15).Value = 18.7
$ws.Cells.Item(26, 16).Value = 24.4
$ws.Cells.Item(26, 17).Value = 0.766
$ws.Cells.Item(26, 19).Value = 28.5
$ws.Cells.Item(26, 21).Value = 20.1
$ws.Cells.Item(26, 25).Value = 3.9
$ws.Cells.Item(26, 26).Value = 20.5
$ws.Cells.Item(26, 28).Value = 98.7
$ws.Cells.Item(26, 29).Value = 3.7
$ws.Cells.Item(26, 30).Value = 8
$ws.Cells.Item(26, 32).Value = 9
$ws.Cells.Item(26, 33).Value = 9
$ws.Cells.Item(26, 34).Value = 14
$ws.Cells.Item(26, 35).Value = 18
$ws.Cells.Item(26, 41).Value = 19
$ws.Cells.Item(26, 43).Value = 17
$ws.Cells.Item(26, 47).Value = 23
$ws.Cells.Item(26, 54).Value = 16
$ws.Cells.Item(26, 55).Value = 7
$ws.Cells.Item(26, 58).Value = "'2009-03-21"
# Row 27
$ws.Cells.Item(27, 30).Value = 8
$ws.Cells.Item(27, 37).Value = 26
$ws.Cells.Item(27, 40).Value = 17
$ws.Cells.Item(27, 45).Value = 26
$ws.Cells.Item(27, 48).Value = 28
$ws.Cells.Item(27, 52).Value = 28
$ws.Cells.Item(27, 58).Value = "'2009-03-21"
# Row 28
$ws.Cells.Item(28, 30).Value = 20
$ws.Cells.Item(28, 55).Value = 6
$ws.Cells.Item(28, 58).Value = "'2009-03-21"
# Row 29
$ws.Cells.Item(29, 30).Value = 8
$ws.Cells.Item(29, 34).Value = 21
$ws.Cells.Item(29, 35).Value = 19
$ws.Cells.Item(29, 58).Value = "'2009-03-21"
# Row 30
$ws.Cells.Item(30, 30).Value = 8
$ws.Cells.Item(30, 31).Value = 8
$ws.Cells.Item(30, 32).Value = 9
$ws.Cells.Item(30, 33).Value = 9
$ws.Cells.Item(30, 35).Value = 6
$ws.Cells.Item(30, 48).Value = 20
$ws.Cells.Item(30, 50).Value = 17
$ws.Cells.Item(30, 58).Value = "'2009-03-21"
# Row 31
$ws.Cells.Item(31, 4).Value = 70
$ws.Cells.Item(31, 6).Value = 54
$ws.Cells.Item(31, 7).Value = 0.229
$ws.Cells.Item(31, 9).Value = 36.3
$ws.Cells.Item(31, 10).Value = 81
$ws.Cells.Item(31, 14).Value = 0.33
$ws.Cells.Item(31, 17).Value = 0.765
$ws.Cells.Item(31, 18).Value = 11.7
$ws.Cells.Item(31, 19).Value = 28
$ws.Cells.Item(31, 20).Value = 39.7
$ws.Cells.Item(31, 22).Value = 14.1
$ws.Cells.Item(31, 25).Value = 5.2
$ws.Cells.Item(31, 27).Value = 19.7
$ws.Cells.Item(31, 28).Value = 95.2
$ws.Cells.Item(31, 29).Value = -7.7
$ws.Cells.Item(31, 30).Value = 2
$ws.Cells.Item(31, 32).Value = 29
$ws.Cells.Item(31, 37).Value = 24
$ws.Cells.Item(31, 41).Value = 26
$ws.Cells.Item(31, 46).Value = 25
$ws.Cells.Item(31, 47).Value = 22
$ws.Cells.Item(31, 51).Value = 21
$ws.Cells.Item(31, 52).Value = 14
$ws.Cells.Item(31, 53).Value = 27
$ws.Cells.Item(31, 58).Value = "'2009-03-21"

Write-Host "Applied all cell updates"